$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "68.150.24"
$ws.Cells.Item(2, 5).Value = "  -0.30%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "3.889.26"
$ws.Cells.Item(3, 5).Value = "  -1.42%  "

# Row 4
$ws.Cells.Item(4, 4).Value = "'0.999"
$ws.Cells.Item(4, 5).Value = "  -0.18%  "

# Row 5
$ws.Cells.Item(5, 4).Value = "'483.01"
$ws.Cells.Item(5, 5).Value = "  +0.16%  "

# Row 6
$ws.Cells.Item(6, 4).Value = "'144.68"
$ws.Cells.Item(6, 5).Value = "  -2.97%  "

# Row 7
$ws.Cells.Item(7, 4).Value = "'0.621"
$ws.Cells.Item(7, 5).Value = "  -0.07%  "

# Row 8
$ws.Cells.Item(8, 4).Value = "'0.998"
$ws.Cells.Item(8, 5).Value = "  -0.01%  "

# Row 9
$ws.Cells.Item(9, 4).Value = "'0.740"
$ws.Cells.Item(9, 5).Value = "  +1.47%  "

# Row 10
$ws.Cells.Item(10, 4).Value = "'0.182"
$ws.Cells.Item(10, 5).Value = "  +9.14%  "

# Row 11
$ws.Cells.Item(11, 4).Value = "'0.0000356"
$ws.Cells.Item(11, 5).Value = "  +1.26%  "

# Row 12
$ws.Cells.Item(12, 4).Value = "'42.87"
$ws.Cells.Item(12, 5).Value = "  +0.07%  "

# Row 13
$ws.Cells.Item(13, 4).Value = "'10.55"
$ws.Cells.Item(13, 5).Value = "  +0.65%  "

# Row 14
$ws.Cells.Item(14, 4).Value = "4.495.00"
$ws.Cells.Item(14, 5).Value = "  -1.48%  "

# Row 15
$ws.Cells.Item(15, 4).Value = "3.890.21"
$ws.Cells.Item(15, 5).Value = "  -2.44%  "

# Row 16
$ws.Cells.Item(16, 4).Value = "'14.27"
$ws.Cells.Item(16, 5).Value = "  -3.27%  "

# Row 17
$ws.Cells.Item(17, 5).Value = "  -0.64%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "'19.99"
$ws.Cells.Item(18, 5).Value = "  +0.27%  "

# Row 19
$ws.Cells.Item(19, 5).Value = "  -0.26%  "

# Row 20
$ws.Cells.Item(20, 4).Value = "68.141.51"
$ws.Cells.Item(20, 5).Value = "  -0.46%  "

# Row 21
$ws.Cells.Item(21, 4).Value = "'429.83"
$ws.Cells.Item(21, 5).Value = "  -1.68%  "

# Row 22
$ws.Cells.Item(22, 5).Value = "  +4.26%  "

# Row 23
$ws.Cells.Item(23, 4).Value = "'14.86"
$ws.Cells.Item(23, 5).Value = "  +2.38%  "

# Row 24
$ws.Cells.Item(24, 4).Value = "'89.80"
$ws.Cells.Item(24, 5).Value = "  +2.57%  "

# Row 25
$ws.Cells.Item(25, 4).Value = "'12.00"
$ws.Cells.Item(25, 5).Value = "  +10.91%  "

# Row 26
$ws.Cells.Item(26, 5).Value = "  +3.13%  "

# Row 27
$ws.Cells.Item(27, 4).Value = "'11.02"
$ws.Cells.Item(27, 5).Value = "  +1.34%  "

# Row 28
$ws.Cells.Item(28, 4).Value = "'37.42"
$ws.Cells.Item(28, 5).Value = "  -2.74%  "

# Row 29
$ws.Cells.Item(29, 5).Value = "  -3.95%  "

# Row 30
$ws.Cells.Item(30, 4).Value = "'710.90"
$ws.Cells.Item(30, 5).Value = "  -1.17%  "

# Row 31
$ws.Cells.Item(31, 4).Value = "'13.54"
$ws.Cells.Item(31, 5).Value = "  +1.76%  "

# Row 32
$ws.Cells.Item(32, 4).Value = "'0.130"
$ws.Cells.Item(32, 5).Value = "  +0.17%  "

# Row 33
$ws.Cells.Item(33, 5).Value = "  +2.03%  "

# Row 34
$ws.Cells.Item(34, 4).Value = "'6.08"
$ws.Cells.Item(34, 5).Value = "  +9.34%  "

# Row 35
$ws.Cells.Item(35, 4).Value = "0.0₃0875"
$ws.Cells.Item(35, 5).Value = "  -2.99%  "

# Row 36
$ws.Cells.Item(36, 2).Value = "InjectiveProtocol"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Cells.Item(36, 4).Value = "'40.96"
$ws.Cells.Item(36, 5).Value = "  -2.81%  "

# Row 37
$ws.Cells.Item(37, 2).Value = "OKB"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Cells.Item(37, 4).Value = "'60.95"
$ws.Cells.Item(37, 5).Value = "  +3.30%  "

# Row 38
$ws.Cells.Item(38, 2).Value = "VeChain"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(38, 4).Value = "'0.0502"
$ws.Cells.Item(38, 5).Value = "  +6.75%  "

# Row 39
$ws.Cells.Item(39, 2).Value = "Kaspa"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(39, 4).Value = "'0.146"
$ws.Cells.Item(39, 5).Value = "  -3.79%  "

# Row 40
$ws.Cells.Item(40, 2).Value = "TheGraph"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Cells.Item(40, 4).Value = "'0.397"
$ws.Cells.Item(40, 5).Value = "  +14.41%  "

# Row 41
$ws.Cells.Item(41, 4).Value = "'0.999"
$ws.Cells.Item(41, 5).Value = "  +0.02%  "

# Row 42
$ws.Cells.Item(42, 5).Value = "  +1.47%  "

# Row 43
$ws.Cells.Item(43, 5).Value = "  +2.60%  "

# Row 44
$ws.Cells.Item(44, 4).Value = "'2.96"
$ws.Cells.Item(44, 5).Value = "  -2.47%  "

# Row 45
$ws.Cells.Item(45, 4).Value = "'0.143"
$ws.Cells.Item(45, 5).Value = "  +0.80%  "

# Row 46
$ws.Cells.Item(46, 5).Value = "  +4.19%  "

# Row 47
$ws.Cells.Item(47, 5).Value = "  -0.11%  "

# Row 48
$ws.Cells.Item(48, 5).Value = "  -1.88%  "

# Row 49
$ws.Cells.Item(49, 4).Value = "'2.12"
$ws.Cells.Item(49, 5).Value = "  -2.87%  "

# Row 50
$ws.Cells.Item(50, 4).Value = "'144.65"
$ws.Cells.Item(50, 5).Value = "  -1.13%  "

# Row 51
$ws.Cells.Item(51, 4).Value = "'2.81"
$ws.Cells.Item(51, 5).Value = "  -1.63%  "
